$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column (D) to be stored as text so numeric-looking
# strings like "1.007" or "325.29" are not auto-converted to numbers,
# matching the inline-string (text) representation used by the feed.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Cells.Item(2, 4).Value = '29.484.52'
$ws.Cells.Item(2, 5).Value = '  +1.29%  '

# Row 3
$ws.Cells.Item(3, 4).Value = '1.914.29'
$ws.Cells.Item(3, 5).Value = '  +1.01%  '

# Row 4
$ws.Cells.Item(4, 4).Value = '1.007'
$ws.Cells.Item(4, 5).Value = '  +0.08%  '

# Row 5
$ws.Cells.Item(5, 4).Value = '325.29'
$ws.Cells.Item(5, 5).Value = '  +1.07%  '

# Row 6
$ws.Cells.Item(6, 5).Value = '  +0.13%  '

# Row 7
$ws.Cells.Item(7, 4).Value = '0.4836'
$ws.Cells.Item(7, 5).Value = '  +3.04%  '

# Row 8
$ws.Cells.Item(8, 4).Value = '0.4073'
$ws.Cells.Item(8, 5).Value = '  +1.57%  '

# Row 9
$ws.Cells.Item(9, 4).Value = '0.08171'
$ws.Cells.Item(9, 5).Value = '  +2.52%  '

# Row 10
$ws.Cells.Item(10, 5).Value = '  +2.82%  '

# Row 11
$ws.Cells.Item(11, 5).Value = '  +5.82%  '

# Row 12
$ws.Cells.Item(12, 4).Value = '1.908.90'
$ws.Cells.Item(12, 5).Value = '  +1.12%  '

# Row 13
$ws.Cells.Item(13, 4).Value = '6.041'
$ws.Cells.Item(13, 5).Value = '  +3.68%  '

# Row 14
$ws.Cells.Item(14, 4).Value = '7.229'
$ws.Cells.Item(14, 5).Value = '  +3.05%  '

# Row 15
$ws.Cells.Item(15, 4).Value = '91.20'
$ws.Cells.Item(15, 5).Value = '  +2.73%  '

# Row 16
$ws.Cells.Item(16, 2).Value = 'TRON'
$ws.Cells.Item(16, 3).Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Cells.Item(16, 4).Value = '0.06775'
$ws.Cells.Item(16, 5).Value = '  +2.33%  '

# Row 17
$ws.Cells.Item(17, 2).Value = 'BinanceUSD'
$ws.Cells.Item(17, 3).Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Cells.Item(17, 4).Value = '1.007'
$ws.Cells.Item(17, 5).Value = '  -0.04%  '

# Row 18
$ws.Cells.Item(18, 4).Value = '0.00001040'
$ws.Cells.Item(18, 5).Value = '  +1.80%  '

# Row 19
$ws.Cells.Item(19, 4).Value = '17.72'
$ws.Cells.Item(19, 5).Value = '  +1.68%  '

# Row 20
$ws.Cells.Item(20, 5).Value = '  +0.54%  '

# Row 21
$ws.Cells.Item(21, 4).Value = '29.510.85'
$ws.Cells.Item(21, 5).Value = '  +1.24%  '

# Row 22
$ws.Cells.Item(22, 4).Value = '5.633'
$ws.Cells.Item(22, 5).Value = '  +2.54%  '

# Row 23
$ws.Cells.Item(23, 5).Value = '  +3.07%  '

# Row 24
$ws.Cells.Item(24, 5).Value = '  -1.07%  '

# Row 25
$ws.Cells.Item(25, 4).Value = '2.159.95'
$ws.Cells.Item(25, 5).Value = '  +2.66%  '

# Row 26
$ws.Cells.Item(26, 4).Value = '6.592'
$ws.Cells.Item(26, 5).Value = '  +10.45%  '

# Row 27
$ws.Cells.Item(27, 4).Value = '156.46'
$ws.Cells.Item(27, 5).Value = '  +1.69%  '

# Row 28
$ws.Cells.Item(28, 4).Value = '20.07'
$ws.Cells.Item(28, 5).Value = '  +2.26%  '

# Row 29
$ws.Cells.Item(29, 4).Value = '2.123'
$ws.Cells.Item(29, 5).Value = '  +2.31%  '

# Row 30
$ws.Cells.Item(30, 4).Value = '120.76'
$ws.Cells.Item(30, 5).Value = '  +3.41%  '

# Row 31
$ws.Cells.Item(31, 4).Value = '1.027'
$ws.Cells.Item(31, 5).Value = '  -2.74%  '

# Row 32
$ws.Cells.Item(32, 4).Value = '0.09542'
$ws.Cells.Item(32, 5).Value = '  +1.12%  '

# Row 33
$ws.Cells.Item(33, 4).Value = '5.517'
$ws.Cells.Item(33, 5).Value = '  +3.62%  '

# Row 34
$ws.Cells.Item(34, 4).Value = '3.556'
$ws.Cells.Item(34, 5).Value = '  -0.17%  '

# Row 35
$ws.Cells.Item(35, 4).Value = '1.394'
$ws.Cells.Item(35, 5).Value = '  -0.06%  '

# Row 36
$ws.Cells.Item(36, 5).Value = '  +2.12%  '

# Row 37
$ws.Cells.Item(37, 4).Value = '0.06127'
$ws.Cells.Item(37, 5).Value = '  +1.30%  '

# Row 38
$ws.Cells.Item(38, 4).Value = '1.189'
$ws.Cells.Item(38, 5).Value = '  +1.94%  '

# Row 39
$ws.Cells.Item(39, 4).Value = '10.89'
$ws.Cells.Item(39, 5).Value = '  +8.94%  '

# Row 40
$ws.Cells.Item(40, 4).Value = '0.5980'
$ws.Cells.Item(40, 5).Value = '  +3.46%  '

# Row 41
$ws.Cells.Item(41, 4).Value = '8.045'
$ws.Cells.Item(41, 5).Value = '  +0.40%  '

# Row 42
$ws.Cells.Item(42, 4).Value = '0.1857'
$ws.Cells.Item(42, 5).Value = '  +2.15%  '

# Row 43
$ws.Cells.Item(43, 4).Value = '2.423'
$ws.Cells.Item(43, 5).Value = '  -1.80%  '

# Row 44
$ws.Cells.Item(44, 5).Value = '  +2.14%  '

# Row 45
$ws.Cells.Item(45, 2).Value = 'Cronos'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Cells.Item(45, 4).Value = '0.07624'
$ws.Cells.Item(45, 5).Value = '  -0.72%  '

# Row 46
$ws.Cells.Item(46, 2).Value = 'EnergySwap'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Cells.Item(46, 4).Value = '12.43'
$ws.Cells.Item(46, 5).Value = '  +2.24%  '

# Row 47
$ws.Cells.Item(47, 4).Value = '0.5596'
$ws.Cells.Item(47, 5).Value = '  +2.63%  '

# Row 48
$ws.Cells.Item(48, 4).Value = '1.956'
$ws.Cells.Item(48, 5).Value = '  +3.43%  '

# Row 49
$ws.Cells.Item(49, 4).Value = '116.58'
$ws.Cells.Item(49, 5).Value = '  +3.27%  '

# Row 50
$ws.Cells.Item(50, 4).Value = '72.84'
$ws.Cells.Item(50, 5).Value = '  +3.07%  '

# Row 51
$ws.Cells.Item(51, 4).Value = '2.424'
$ws.Cells.Item(51, 5).Value = '  +3.62%  '
